$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.914.62"
$ws.Range("E2").Value = "  +1.60%  "
$ws.Range("D3").Value = "2.125.62"
$ws.Range("E3").Value = "  +11.04%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "257.70"
$ws.Range("E5").Value = "  +3.73%  "
$ws.Range("D6").Value = "0.674"
$ws.Range("E6").Value = "  -2.69%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "47.36"
$ws.Range("E8").Value = "  +8.07%  "
$ws.Range("D9").Value = "61.63"
$ws.Range("E9").Value = "  +6.83%  "
$ws.Range("D10").Value = "0.377"
$ws.Range("E10").Value = "  +3.46%  "
$ws.Range("D11").Value = "0.0746"
$ws.Range("E11").Value = "  -1.95%  "
$ws.Range("E12").Value = "  +1.19%  "
$ws.Range("D13").Value = "2.429.01"
$ws.Range("E13").Value = "  +10.55%  "
$ws.Range("E14").Value = "  +1.78%  "
$ws.Range("E15").Value = "  +8.12%  "
$ws.Range("D16").Value = "2.120.82"
$ws.Range("E16").Value = "  +10.48%  "
$ws.Range("E17").Value = "  +1.69%  "
$ws.Range("D18").Value = "36.809.62"
$ws.Range("E18").Value = "  +1.35%  "
$ws.Range("D19").Value = "74.52"
$ws.Range("E19").Value = "  +0.41%  "
$ws.Range("D20").Value = "0.0₃0850"
$ws.Range("E20").Value = "  +0.52%  "
$ws.Range("D21").Value = "13.56"
$ws.Range("E21").Value = "  +3.10%  "
$ws.Range("D22").Value = "243.59"
$ws.Range("E22").Value = "  -3.21%  "
$ws.Range("D23").Value = "5.26"
$ws.Range("E23").Value = "  +1.88%  "
$ws.Range("E24").Value = "  +0.12%  "
$ws.Range("D25").Value = "2.47"
$ws.Range("E25").Value = "  -8.35%  "
$ws.Range("D26").Value = "173.15"
$ws.Range("E26").Value = "  +3.31%  "
$ws.Range("D27").Value = "21.59"
$ws.Range("E27").Value = "  +15.04%  "
$ws.Range("D28").Value = "9.31"
$ws.Range("E28").Value = "  +6.22%  "
$ws.Range("E29").Value = "  -6.11%  "
$ws.Range("E30").Value = "  -3.28%  "
$ws.Range("D31").Value = "23.42"
$ws.Range("E31").Value = "  +52.63%  "
$ws.Range("E32").Value = "  +1.37%  "
$ws.Range("D33").Value = "0.0961"
$ws.Range("E33").Value = "  +14.07%  "
$ws.Range("D34").Value = "0.0606"
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("D35").Value = "2.45"
$ws.Range("E35").Value = "  +22.46%  "
$ws.Range("D36").Value = "4.25"
$ws.Range("E36").Value = "  -1.42%  "
$ws.Range("E37").Value = "  -0.10%  "
$ws.Range("D38").Value = "1.89"
$ws.Range("E38").Value = "  -3.41%  "
$ws.Range("D39").Value = "0.927"
$ws.Range("E39").Value = "  +8.23%  "
$ws.Range("D40").Value = "1.38"
$ws.Range("E40").Value = "  -6.90%  "
$ws.Range("D41").Value = "1.20"
$ws.Range("E41").Value = "  +9.07%  "
$ws.Range("D42").Value = "0.0226"
$ws.Range("E42").Value = "  -0.97%  "
$ws.Range("D43").Value = "100.17"
$ws.Range("E43").Value = "  -4.88%  "
$ws.Range("D44").Value = "2.82"
$ws.Range("E44").Value = "  +16.45%  "
$ws.Range("D45").Value = "16.40"
$ws.Range("E45").Value = "  -3.63%  "
$ws.Range("D46").Value = "1.367.51"
$ws.Range("E46").Value = "  +2.14%  "
$ws.Range("D47").Value = "0.0843"
$ws.Range("E47").Value = "  +4.37%  "
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").Value = "7.22"
$ws.Range("E48").Value = "  +12.46%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "2.33"
$ws.Range("E49").Value = "  -1.72%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.306.00"
$ws.Range("E50").Value = "  +9.98%  "
$ws.Range("D51").Value = "2.83"
$ws.Range("E51").Value = "  +1.63%  "
